$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Results for match 48 (SRH vs DC) - row 56
$ws.Range("E56").Value = 60
$ws.Range("H56").Value = 80
$ws.Range("K56").Value = 40
$ws.Range("N56").Value = 100
$ws.Range("Q56").Value = 20
$ws.Range("T56").Value = 0

# Results for match 49 (MI vs RCB) - row 57
$ws.Range("E57").Value = 20
$ws.Range("H57").Value = 60
$ws.Range("K57").Value = 100
$ws.Range("N57").Value = 40
$ws.Range("Q57").Value = 80
$ws.Range("T57").Value = 0
